# "Rerun the example data" - update CoreSNPclusters_ST39_SNPcutoff20_Days45
#
# The underlying analysis was rerun, which:
#  - changed cluster/group id "3" -> "2" (affects existing rows that were in
#    group "3": row 3 and row 4)
#  - changed row 4 (id=3)'s isolate label from EXM2144 -> EXM2170, with an
#    updated collection date (2020-09-13 instead of 2020-09-30) and epiweek
#    38 instead of 40 (epiyearweek 2020.38)
#  - added a new row (id=4) for isolate EXM2144 using the data that row 4
#    used to hold (2020-09-30, epiweek 40, epiyearweek 2020.40), now in
#    group "2" as well

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- update existing row 3 (table data row 2 / sheet row 3): group 3 -> 2 ---
$ws.Cells.Item(3, 10).Value = "2"

# --- update existing row 4 (table data row 3 / sheet row 4) ---
$ws.Cells.Item(4, 2).Value = "EXM2170"     # label
$ws.Cells.Item(4, 3).Value = 44087          # Collectiondate serial
$ws.Cells.Item(4, 8).Value = 38             # epiweek
$ws.Cells.Item(4, 9).Value = "2020.38"      # epiyearweek
$ws.Cells.Item(4, 10).Value = "2"           # group

# --- add new row 5 (table data row 4) for id 4 / EXM2144 ---
$newRow = $tbl.ListRows.Add()
$newRange = $newRow.Range

# copy the date cell's number formatting from the row above so the new
# date cell keeps a consistent, shared style
$ws.Range("C4").Copy()
$newRange.Cells.Item(1, 3).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRange.Cells.Item(1, 1).Value = 4                       # id
$newRange.Cells.Item(1, 2).Value = "EXM2144"                # label
$newRange.Cells.Item(1, 3).Value = 44104                    # Collectiondate
$newRange.Cells.Item(1, 4).Value = "EXAMPLE HOSPITAL"        # FacilityName
$newRange.Cells.Item(1, 5).Value = "neonatal"                # WardType
$newRange.Cells.Item(1, 6).Value = 39                        # ST
$newRange.Cells.Item(1, 7).Value = 2020                      # epiyear
$newRange.Cells.Item(1, 8).Value = 40                        # epiweek
$newRange.Cells.Item(1, 9).Value = "2020.40"                 # epiyearweek
$newRange.Cells.Item(1, 10).Value = "2"                      # group
$newRange.Cells.Item(1, 11).Value = "UNKNOWN"                # X1
$newRange.Cells.Item(1, 12).Value = "UNKNOWN"                # X2
$newRange.Cells.Item(1, 13).Value = "UNKNOWN"                # X3

# --- apply the (now builtin) short-date format to the whole date column ---
$ws.Range("C2:C5").NumberFormat = "mm-dd-yy"
